# CodeSystem-nmdp-abo-group-cs.xlsx : "Update for release to deploy 0.1.1"
#
# Changes to the "Metadata" sheet:
#   - Version (B3): 0.1.0 -> 0.1.1
#   - Date    (B8): 2024-05-20T10:22:59-05:00 -> 2024-11-11T17:53:38-06:00
#   - A new "Jurisdiction" property row (with an empty Value) is inserted
#     right after the existing "Contact" row, pushing every row below it
#     down by one (old row 11 "Description" ... old row 21 "Count" become
#     new rows 12 ... 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update Version and Date values -----------------------------------
$ws.Cells.Item(3, 2).Value2 = "0.1.1"
$ws.Cells.Item(8, 2).Value2 = "2024-11-11T17:53:38-06:00"

# --- Insert the new "Jurisdiction" row at row 11 -----------------------
$ws.Rows.Item(11).Insert()

# Give the freshly-inserted row the same formatting as the rest of the
# property/value table (row-insert otherwise invents a slightly different
# style), by copying format from the row right below it.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value2 = "Jurisdiction"

# Leading apostrophe forces a genuine (shared-string) empty-text value in
# B11 instead of Excel just clearing the cell to blank.
$ws.Cells.Item(11, 2).Value = "'"

# The apostrophe trick stamps a "quote prefix" onto the cell's style;
# reapply the plain table formatting so B11 ends up identical in style to
# its neighbours.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$excel.CutCopyMode = $false
